$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of the existing header cell (H1) onto the two new header
# cells so they pick up the same bold/bordered/centered formatting (style index 1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data column I: constant value of 1 for every data row
$ws.Range("I2").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("I7").Value = 1

# New data column J: mirrors column H's values
$ws.Range("J2").Value = 4
$ws.Range("J3").Value = 5
$ws.Range("J4").Value = 4
$ws.Range("J5").Value = 3
$ws.Range("J6").Value = 3
$ws.Range("J7").Value = 2
